# Update capital structure database values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "X2"  = 0.06119801074845367
    "Y2"  = 0.04012212126474764
    "AB2" = 0.05889485741528916
    "AC2" = -0.05889485741528916

    "X3"  = 0.03258487127035474
    "Y3"  = 0.03005663816360753
    "AB3" = 0.04137483310039668
    "AC3" = 0.02798378819199775

    "X4"  = 0.03606982588866234
    "Y4"  = 0.1011458701873187
    "AB4" = 0.04727265882077926
    "AC4" = -0.04727265882077926

    "X5"  = 0.03115429517334863
    "Y5"  = 0.06852621600875999
    "AB5" = 0.0504218527544029
    "AC5" = -0.0504218527544029

    "X6"  = 0.04897825078191798
    "Y6"  = 0.07639488354644022
    "AB6" = 0.05150557567097748
    "AC6" = -0.05150557567097748

    "X7"  = 0.06119801074845367
    "Y7"  = 0.04012212126474764
    "AB7" = 0.05889485741528916
    "AC7" = -0.05889485741528916

    "X8"  = 0.0703741747904381
    "Y8"  = 0.0006564658780855714
    "AB8" = 0.06081726831145565
    "AC8" = -0.06081726831145565

    "X9"  = 0.07265592769655503
    "Y9"  = 0.1116788362519428
    "AB9" = 0.06120359006866481
    "AC9" = -0.06120359006866481

    "X10"  = 0.1427384776344401
    "Y10"  = -0.03967441078207237
    "AB10" = 0.06645445645892989
    "AC10" = -0.06645445645892989

    "X11"  = 0.06310970120336465
    "Y11"  = -0.007591306554535215
    "AB11" = 0.05935333823834295
    "AC11" = -0.07335812025627025
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
